$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6748145818710327
$ws.Range("B1").Value = 0.7091183662414551
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.528087615966797
$ws.Range("E1").Value = 0.913688063621521
